# Jenkins.sh build run: switch the automation run's "environment" parameter
# from a local run to a remote (SauceLabs) run, and leave the Parameters
# sheet active/selected at the updated cell, matching how this was edited
# and saved in Excel.

$wb = $excel.ActiveWorkbook

$parameters = $wb.Worksheets.Item("Parameters")

# Update the environment value used by the run: "local" -> "remote"
$parameters.Range("B2").Value = "remote"

# Make "Parameters" the active sheet/tab, with B2 selected - this is the
# sheet+cell that was active when the workbook was saved after the edit.
$parameters.Activate() | Out-Null
$parameters.Range("B2").Select() | Out-Null
